$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

$msgObjRef = "AppEx: Object reference not set to an instance of an object. at Source: Invoke ACME_OpenAndLogin workflow"
$msgFile = "AppEx: Could not find file 'C:\Users\smit.rathore\Documents\UiPath\WI4-Performer-Demo-Can\Test_Framework\Test_ProcessTransaction.xaml'. at Source: mscorlib"
$msgSelector = "AppEx: Cannot find the UI element corresponding to this selector: <html title='ACME System 1*' /> at Source: Invoke ACME_LogoutAndClose workflow: Attach Browser 'iexplore.exe ACME'"

$rows = @(
    @{Row = 3;  Comment = $msgObjRef},
    @{Row = 8;  Comment = $msgObjRef},
    @{Row = 9;  Comment = $msgFile},
    @{Row = 10; Comment = $msgFile},
    @{Row = 11; Comment = $msgFile},
    @{Row = 12; Comment = $msgFile},
    @{Row = 13; Comment = $msgFile},
    @{Row = 14; Comment = $msgFile},
    @{Row = 15; Comment = $msgFile},
    @{Row = 16; Comment = $msgFile},
    @{Row = 17; Comment = $msgFile},
    @{Row = 18; Comment = $msgSelector}
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = "FAIL"
    $ws.Cells.Item($r, 4).Value = $entry.Comment
}
